$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.2
$ws.Range("H3").Value = 2.4
$ws.Range("K3").Value = 5
$ws.Range("G4").Value = 2.98
$ws.Range("J4").Value = 2.68
$ws.Range("P4").Value = 1.49
$ws.Range("Q4").Value = 2.64
$ws.Range("F7").Value = 2.58
$ws.Range("G7").Value = 2.6
$ws.Range("H7").Value = 2.78
$ws.Range("I7").Value = 2.86
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 3.95
$ws.Range("H8").Value = 2.02
$ws.Range("I8").Value = 2.16
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 3.7
$ws.Range("P8").Value = 1.76
$ws.Range("P11").Value = 1.78
$ws.Range("G12").Value = 1.75
$ws.Range("H12").Value = 6.2
$ws.Range("I12").Value = 7.4
$ws.Range("K12").Value = 4.2
$ws.Range("P12").Value = 1.69
$ws.Range("Q12").Value = 2.16
$ws.Range("G13").Value = 2.16
$ws.Range("J13").Value = 3.25
$ws.Range("K13").Value = 5.3
$ws.Range("F14").Value = 1.91
$ws.Range("G14").Value = 2.02
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4.6
$ws.Range("J14").Value = 3.6
$ws.Range("K14").Value = 4.1
$ws.Range("Q14").Value = 1.84
$ws.Range("F15").Value = 3.9
$ws.Range("G15").Value = 3.95
$ws.Range("H15").Value = 2.14
$ws.Range("I15").Value = 2.18
$ws.Range("V15").Value = 1.84
$ws.Range("W15").Value = 1.33
$ws.Range("Z15").Value = 13
$ws.Range("AE15").Value = 24
$ws.Range("P16").Value = 1.81
$ws.Range("G17").Value = 2.06
$ws.Range("I17").Value = 5.8
$ws.Range("K17").Value = 3.3
